# Luke Dysart final push
# Reorders the displayed text of the hyperlinked part-name cells in row 3
# (F3, G3, H3) so that:
#   F3 (was GFP)    -> B0034
#   G3 (was B0034)  -> B0015
#   H3 (was B0015)  -> GFP

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "https://sbolcanvas.org/B0034/1"
$ws.Range("G3").Value = "https://sbolcanvas.org/B0015/1"
$ws.Range("H3").Value = "https://sbolcanvas.org/GFP/1"
